# 20170808 功能微調 bug fixed
#
# Changes applied (per the OOXML diff):
#  1. Add a workbook-level print area defined name for the one sheet:
#     'Print_Area' = '完整版-航空公司'!$A$1:$G$96
#  2. Tighten the page margins to the "narrow" metric set (0.6/1.9/0.8 cm)
#     expressed in points, matching 0.23622047244094491in /
#     0.74803149606299213in / 0.31496062992125984in.
#  3. Shrink the print scale to 88%.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Print area -----------------------------------------------------------
$ws.PageSetup.PrintArea = '$A$1:$G$96'

# -- Page margins (points = inches * 72) -----------------------------------
$ws.PageSetup.LeftMargin   = 17.007874015748033   # 0.23622047244094491 in
$ws.PageSetup.RightMargin  = 17.007874015748033   # 0.23622047244094491 in
$ws.PageSetup.TopMargin    = 53.85826771653544    # 0.74803149606299213 in
$ws.PageSetup.BottomMargin = 53.85826771653544    # 0.74803149606299213 in
$ws.PageSetup.HeaderMargin = 22.677165354330707   # 0.31496062992125984 in
$ws.PageSetup.FooterMargin = 22.677165354330707   # 0.31496062992125984 in

# -- Print scale ------------------------------------------------------------
$ws.PageSetup.Zoom = 88
